$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (organiser_id and everything to the right shifts
# from C..J to D..K).
$ws.Columns.Item(3).Insert()

# New "category" header + per-row category values.
$ws.Range("C1").Value = "category"
$ws.Range("C2").Value = "other"
$ws.Range("C3").Value = "other"
$ws.Range("C4").Value = "sport"
$ws.Range("C5").Value = "food and drink"
$ws.Range("C6").Value = "music"
$ws.Range("C7").Value = "art"
$ws.Range("C8").Value = "music"
$ws.Range("C9").Value = "family"
$ws.Range("C10").Value = "sport"
$ws.Range("C11").Value = "sport"

# Widen the new column (it does not inherit the bestFit auto-width of its
# neighbours).
$ws.Columns.Item(3).ColumnWidth = 35

# Match the author's final selection/cursor position.
$ws.Range("C12").Select() | Out-Null
